$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value2 = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "97.042.63"
$ws.Range("E2").Value2 = "  +2.15%  "
Set-TextValue "D3" "3.591.25"
$ws.Range("E3").Value2 = "  -0.67%  "
$ws.Range("E4").Value2 = "  +0.04%  "
Set-TextValue "D5" "243.60"
$ws.Range("E5").Value2 = "  +3.29%  "
Set-TextValue "D6" "654.71"
$ws.Range("E6").Value2 = "  -0.21%  "
Set-TextValue "D7" "1.65"
$ws.Range("E7").Value2 = "  +13.90%  "
Set-TextValue "D8" "0.413"
$ws.Range("E8").Value2 = "  +3.65%  "
Set-TextValue "D9" "1.06"
$ws.Range("E9").Value2 = "  +7.02%  "
Set-TextValue "D10" "1.00"
$ws.Range("E10").Value2 = "  +0.01%  "
Set-TextValue "D11" "3.587.96"
$ws.Range("E11").Value2 = "  -0.75%  "
Set-TextValue "D12" "43.72"
$ws.Range("E12").Value2 = "  +4.15%  "
$ws.Range("E13").Value2 = "  +1.51%  "
Set-TextValue "D14" "6.45"
$ws.Range("E14").Value2 = "  +0.57%  "
Set-TextValue "D15" "4.256.34"
$ws.Range("E15").Value2 = "  -1.00%  "
Set-TextValue "D16" "96.772.02"
$ws.Range("E16").Value2 = "  +1.94%  "
$ws.Range("E17").Value2 = "  +2.29%  "
Set-TextValue "D18" "3.583.19"
$ws.Range("E18").Value2 = "  -0.70%  "
Set-TextValue "D19" "7.76"
$ws.Range("E19").Value2 = "  -1.77%  "
Set-TextValue "D20" "12.73"
$ws.Range("E20").Value2 = "  -0.95%  "
Set-TextValue "D21" "18.05"
$ws.Range("E21").Value2 = "  +0.65%  "
Set-TextValue "D22" "0.531"
$ws.Range("E22").Value2 = "  +11.45%  "
Set-TextValue "D23" "509.59"
$ws.Range("E23").Value2 = "  +1.66%  "
Set-TextValue "D24" "3.42"
$ws.Range("E24").Value2 = "  -3.01%  "
Set-TextValue "D25" "0.0000202"
$ws.Range("E25").Value2 = "  +3.62%  "
Set-TextValue "D26" "6.93"
$ws.Range("E26").Value2 = "  +5.59%  "
$ws.Range("E27").Value2 = "  +2.79%  "
Set-TextValue "D28" "13.18"
$ws.Range("E28").Value2 = "  +6.07%  "
Set-TextValue "D29" "3.781.26"
$ws.Range("E29").Value2 = "  -0.65%  "
$ws.Range("B30").Value2 = "PancakeSwap"
$ws.Range("C30").Value2 = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D30" "3.05"
$ws.Range("E30").Value2 = "  -1.32%  "
$ws.Range("B31").Value2 = "Hedera"
$ws.Range("C31").Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D31" "0.152"
$ws.Range("E31").Value2 = "  +10.79%  "
Set-TextValue "D32" "11.53"
$ws.Range("E32").Value2 = "  +3.03%  "
Set-TextValue "D33" "0.999"
$ws.Range("E33").Value2 = "  -0.20%  "
$ws.Range("E34").Value2 = "  +6.12%  "
Set-TextValue "D35" "0.999"
$ws.Range("E35").Value2 = "  -0.21%  "
Set-TextValue "D36" "31.52"
$ws.Range("E36").Value2 = "  -2.05%  "
Set-TextValue "D37" "8.96"
$ws.Range("E37").Value2 = "  +11.97%  "
Set-TextValue "D38" "627.73"
$ws.Range("E38").Value2 = "  +11.56%  "
$ws.Range("E39").Value2 = "  +2.92%  "
$ws.Range("E40").Value2 = "  +12.67%  "
$ws.Range("E41").Value2 = "  +1.62%  "
$ws.Range("E42").Value2 = "  +0.04%  "
Set-TextValue "D43" "0.915"
$ws.Range("E43").Value2 = "  +0.43%  "
$ws.Range("E44").Value2 = "  +7.07%  "
Set-TextValue "D45" "5.83"
$ws.Range("E45").Value2 = "  +3.67%  "
$ws.Range("E46").Value2 = "  +5.66%  "
$ws.Range("E47").Value2 = "  +3.51%  "
Set-TextValue "D48" "23.56"
$ws.Range("E48").Value2 = "  -0.50%  "
Set-TextValue "D49" "33.29"
$ws.Range("E49").Value2 = "  -7.89%  "
$ws.Range("E50").Value2 = "  +0.91%  "
Set-TextValue "D51" "8.31"
$ws.Range("E51").Value2 = "  +4.17%  "
